$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample rows: tenant "government" with two service paths, plus a
# NULL-tenant / NULL-path row mirroring the existing row 2 pattern.
$ws.Range("A4").Value = "NULL"
$ws.Range("B4").Value = "government"

$ws.Range("A5").Value = "/park"
$ws.Range("B5").Value = "government"

$ws.Range("A6").Value = "/inbound"
$ws.Range("B6").Value = "government"
